# Rename existing sheet and add the new "LCV" sheet right after it.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "binek_arac"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "LCV"

# --- Shared strings / labels on binek_arac -------------------------------
# Row 2 label changes from "toplam_hurda_tesviki_butcesi" to
# "yeni_arac_indirim_orani_grup1"; row 3 label changes from
# "yeni_arac_indirim_orani" to "yeni_arac_indirim_orani_grup2".
$ws1.Range("A2").Value = "yeni_arac_indirim_orani_grup1"
$ws1.Range("A3").Value = "yeni_arac_indirim_orani_grup2"

# Header row is bold.
$ws1.Range("A1:B1").Font.Bold = $true

# Values: B2 becomes 0.2 (rate instead of budget amount), B3 stays 0.15
# but loses the old "Virgul" (thousands) number style.
$ws1.Range("B2").ClearFormats()
$ws1.Range("B2").Value = 0.2
$ws1.Range("B3").Value = 0.15

# Selection / active cell on binek_arac after the edits.
[void]$ws1.Range("A2:B2").Select()

# Print setup (A4 portrait) picked up by the sheet on save.
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# --- New "LCV" sheet -------------------------------------------------------
$ws2.Cells.Item(1, 1).Value = "degisken"
$ws2.Cells.Item(1, 2).Value = "deger"
$ws2.Range("A1:B1").Font.Bold = $true

$ws2.Cells.Item(2, 1).Value = "LCV_hurda_tesvik_orani"
$ws2.Cells.Item(2, 2).Value = 0.15

$ws2.Columns.Item(1).ColumnWidth = 28.1666666666667

# Make LCV the active (selected) sheet/tab, matching activeTab="1".
[void]$ws2.Activate()
[void]$ws2.Range("I10").Select()
